$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of k (column J) across the 10 instances
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Summary statistics block (rows 14-17)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Bold, 12pt, vertically-centered style for the summary values; apply to
# B14 first then copy the resulting format onto B15:B17 so only a single
# new style entry is interned (mirrors how the source workbook resolves
# to one extra cellXf).
$r = $ws.Range("B14")
$r.Font.Bold = $true
$r.Font.Size = 12
$r.VerticalAlignment = -4108

$r.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# The taller font bumps the row height for the summary block.
$ws.Rows("14:17").RowHeight = 15.6

# Page setup as left by the author's last save.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Final selection left active by the author.
$ws.Range("A14:B17").Select()
